$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line: remove "EF (9), EM (8), EB (8), "
$d.Content.Find.Execute(
    "Curso (semestre ideal): EF (9), EM (8), EB (8), EP (10), EQD (8), EQN (11)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Curso (semestre ideal): EP (10), EQD (8), EQN (11)", 2)

# 2. Remove the "Requisitos" heading paragraph and the following list-bullet
#    paragraph ("LOB1008 - ... (Requisito fraco)") at the end of the document.
$paraCount = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Requisitos") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($targetIndex)
    $endPara = $d.Paragraphs.Item($targetIndex + 1)
    $start = $startPara.Range.Start
    $end = $endPara.Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
